# Add two new JETT expression cells to the "Second" worksheet, exposing
# the POI "cell" object (row/column index + cell style wrap-text state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Second")

# New cell at B5: reports the implicit "cell" object's row/column index.
$ws.Range("B5").Value = 'This Cell is at row ${cell.rowIndex} and column ${cell.columnIndex}.'

# New cell at D7: reports whether the implicit "cell" object wraps text,
# and is itself styled with wrap text enabled (creates a new cellXfs entry).
$ws.Range("D7").Value = 'This Cell''s text is ${cell.cellStyle.wrapText ? "wrapped." : "not wrapped."}'
$ws.Range("D7").WrapText = $true

# Make the new row/column visible: widen column D and give row 7 extra
# height so the wrapped text has room to display.
$ws.Rows.Item(7).RowHeight = 60
$ws.Columns.Item(4).ColumnWidth = 17
